# shebang added for executable support
#
# - Rename the "Chips" header (C1) to "Snack"
# - Log a new snack entry as row 8: same date as row 7 (2018-11-20 /
#   serial 43424), time 12:55 PM, snack "Pocky (Chocolate)"
# - Merge the new row's date cell with row 7's date cell (A7:A8), the
#   same way the existing A2:A5 block groups same-day entries
# - Move the active selection to the newly entered time cell (B8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header rename: "Chips" -> "Snack" ---------------------------------
$ws.Range("C1").Value = "Snack"

# --- populate the new row (8) -------------------------------------------
$ws.Range("A8").Value = $ws.Range("A7").Value2
$ws.Range("B8").Value = 0.53819444444444442
$ws.Range("C8").Value = "Pocky (Chocolate)"

# carry over row 7's date/time formatting onto row 8 via copy/paste-format
# so the existing style entries are reused instead of new ones invented
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# --- group the same-day date cells, like the existing A2:A5 merge ------
$ws.Range("A7:A8").Merge()

# --- move the selection to the newly added time cell ---------------------
$ws.Range("B8").Select()
